$d = $word.ActiveDocument
$nl = [char]11
$warnings = 0

# --- Paragraph 6 ---
$p = $d.Paragraphs.Item(6)
$old = "Fornecer ao estudante os principais tipos de síntese orgânica e inorgânica de materiais bem como apresentar as principais técnicas analíticas para caracterização de materiais."
$new = "Introdução à química e sua associação com síntese de novos materiais. A visão moderna do átomo  e Ligações químicas. Estrutura cristalina e técnicas de caracterização cristalográfica. Filmes finos epitaxiais e filmes de uma maneira geral e seu impacto na tecnologica moderna. Crescimento de cristais  Materiais amorfos, síntese e aplicações. Processos e Técnicas de crescimento de cristais de um modo geral. Polímeros condutores e suas aplicações em tecnologica moderna."
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 6 run 0 not found"; $warnings++ }

# --- Paragraph 7 ---
$p = $d.Paragraphs.Item(7)
$old = "Provide the student with the main types of organic and inorganic synthesis of materials as well as presenting the main analytical techniques for material characterization."
$new = "Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology."
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 7 run 0 not found"; $warnings++ }

# --- Paragraph 9 ---
$p = $d.Paragraphs.Item(9)
$old = "5840730 - Antonio Jefferson da Silva Machado" + $nl + ""
$new = "Fornecer ao estudante os principais tipos de síntese orgânica e inorgânica de materiais bem como apresentar as principais técnicas analíticas para caracterização de materiais." + $nl + ""
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 9 run 0 not found"; $warnings++ }
$old = "5840897 - Clodoaldo Saron"
$new = "Química de materiais: definição; papel da química na ciência de materiais; fundamentos." + $nl + "Atomística e a visão moderna do átomo com fundamentos quânticos.Tipos de ligações químicas: forças de van der Waals, potencial de Lennard-Jones, ligação covalente, ligações por coordenação, ligações iônicas e ligações metálicas." + $nl + "Materiais policristalinos e monocristalinos. A ordem cristalográfica e técnicas de caracterização cristalográfica e microscópica. A importância de monocristais em aplicações eletrônicas. Técnicas de crescimento de cristais de alta qualidade tais como: método do fluxo, método Czochralski, método Brigdmann, método do transporte de vapor e método de crescimento de transporte de vapor modificado e isotérmico. Materiais amorfos e sua importância para a tecnologica moderna. Conceitos e técnicas de crescimento de materiais amorfos. Filmes finos epitaxiais, técnicas de crescimento tais como: vapor químico, sputtering, laser ablation e MBE. Filmes finos crescidos por eletrólise para revestimento protetivo, conceitos e aplicações. Síntese de polímeros condutores, conceitos e aplicações como dispositivos eletrônicos."
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 9 run 1 not found"; $warnings++ }

# --- Paragraph 11 ---
$p = $d.Paragraphs.Item(11)
$old = "Introdução à química e sua associação com síntese de novos materiais. A visão moderna do átomo  e Ligações químicas. Estrutura cristalina e técnicas de caracterização cristalográfica. Filmes finos epitaxiais e filmes de uma maneira geral e seu impacto na tecnologica moderna. Crescimento de cristais  Materiais amorfos, síntese e aplicações. Processos e Técnicas de crescimento de cristais de um modo geral. Polímeros condutores e suas aplicações em tecnologica moderna."
$new = "Aulas expositivas e práticas ministradas em laboratório."
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 11 run 0 not found"; $warnings++ }

# --- Paragraph 12 ---
$p = $d.Paragraphs.Item(12)
$old = "Introduction to the chemistry of materials and its association with the synthesis of new materials. The modern view of the atom and chemical bonds. Crystal structure and crystallographic characterization techniques. Epitaxial thin films and films in general and their impact on modern technology. Amorphous materials, synthesis and applications. Synthesis of materials and chemical transformations. Processes and Techniques of crystal growth in general. Conducting polymers and their applications in modern technology."
$new = "Provide the student with the main types of organic and inorganic synthesis of materials as well as presenting the main analytical techniques for material characterization."
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 12 run 0 not found"; $warnings++ }

# --- Paragraph 14 ---
$p = $d.Paragraphs.Item(14)
$old = "Química de materiais: definição; papel da química na ciência de materiais; fundamentos." + $nl + "Atomística e a visão moderna do átomo com fundamentos quânticos.Tipos de ligações químicas: forças de van der Waals, potencial de Lennard-Jones, ligação covalente, ligações por coordenação, ligações iônicas e ligações metálicas." + $nl + "Materiais policristalinos e monocristalinos. A ordem cristalográfica e técnicas de caracterização cristalográfica e microscópica. A importância de monocristais em aplicações eletrônicas. Técnicas de crescimento de cristais de alta qualidade tais como: método do fluxo, método Czochralski, método Brigdmann, método do transporte de vapor e método de crescimento de transporte de vapor modificado e isotérmico. Materiais amorfos e sua importância para a tecnologica moderna. Conceitos e técnicas de crescimento de materiais amorfos. Filmes finos epitaxiais, técnicas de crescimento tais como: vapor químico, sputtering, laser ablation e MBE. Filmes finos crescidos por eletrólise para revestimento protetivo, conceitos e aplicações. Síntese de polímeros condutores, conceitos e aplicações como dispositivos eletrônicos."
$new = "Média simples de duas provas escritas,  Conceito Final = (P1 + P2)/2"
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 14 run 0 not found"; $warnings++ }

# --- Paragraph 17 ---
$p = $d.Paragraphs.Item(17)
$old = "Aulas expositivas e práticas ministradas em laboratório." + $nl + ""
$new = "Aplicação de duas provas escritas dentro do prazo regimental antes do início do próximo semestre letivo." + $nl + ""
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 17 run 1 not found"; $warnings++ }
$old = "Média simples de duas provas escritas,  Conceito Final = (P1 + P2)/2" + $nl + ""
$new = "ALLCOCK, H. R. Introduction to Materials Chemistry, Wiley, Nova Iorque, 2008." + $nl + "FAHLMAN, B. D. Materials Chemistry, Springer, Holanda, 2007." + $nl + "ZHANG, S.; LI, L.; KUMAR, A. Materials Characterization Techniques, Boca Raton: CRC Press, 2008." + $nl + "LENG, Y. Materials Characterization: Introduction to Microscopic and Spectroscopic Methods, Wiley, Cingapura, 2008." + $nl + ""
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 17 run 3 not found"; $warnings++ }
$old = "Aplicação de duas provas escritas dentro do prazo regimental antes do início do próximo semestre letivo."
$new = "5840730 - Antonio Jefferson da Silva Machado"
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 17 run 5 not found"; $warnings++ }

# --- Paragraph 19 ---
$p = $d.Paragraphs.Item(19)
$old = "ALLCOCK, H. R. Introduction to Materials Chemistry, Wiley, Nova Iorque, 2008." + $nl + "FAHLMAN, B. D. Materials Chemistry, Springer, Holanda, 2007." + $nl + "ZHANG, S.; LI, L.; KUMAR, A. Materials Characterization Techniques, Boca Raton: CRC Press, 2008." + $nl + "LENG, Y. Materials Characterization: Introduction to Microscopic and Spectroscopic Methods, Wiley, Cingapura, 2008."
$new = "5840897 - Clodoaldo Saron"
$found = $p.Range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $found) { Write-Host "WARN: paragraph 19 run 0 not found"; $warnings++ }

Write-Host "Done. Warnings:" $warnings